$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the previous single-column layout entirely before laying out the
# new two-column "first name / last name" table.
$ws.Range("A1:B17").ClearContents()

# Column headers
$ws.Range("A1").Value = "first name"
$ws.Range("B1").Value = "last name"

# Row 2, column A keeps its original `="a"` formula.
$ws.Range("A2").Formula = '="a"'
$ws.Range("B2").Value = "aa"

$ws.Range("A3").Value = "b"
$ws.Range("B3").Value = "bb"

$ws.Range("A4").Value = "c"
$ws.Range("B4").Value = "cc"

$ws.Range("A5").Value = "d"
$ws.Range("B5").Value = "dd"

$ws.Range("A7").Value = "f"
$ws.Range("B7").Value = "ff"

$ws.Range("A8").Value = "g"
$ws.Range("B8").Value = "gg"

$ws.Range("A9").Value = "h"
$ws.Range("B9").Value = "hh"

$ws.Range("A6").Value = "er"
$ws.Range("B6").Value = "ee"

$ws.Range("A10").Value = "y"
$ws.Range("B10").Value = "ii"

$ws.Range("A11").Value = "t"
$ws.Range("B11").Value = "jj"

$ws.Range("A12").Value = "jh"
$ws.Range("B12").Value = "dg"

$ws.Range("A13").Value = "f"
$ws.Range("B13").Value = "dfg"

$ws.Range("A14").Value = "d"
$ws.Range("B14").Value = "dfg"

$ws.Range("A15").Value = "lii"
$ws.Range("B15").Value = "dfg"

$ws.Range("A16").Value = "egr"
$ws.Range("B16").Value = "dfg"

# Match the final selection/window state recorded in the saved file.
$ws.Range("G11").Select()
